$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Fix the typo "instantanné" -> "instantané" in the messagerie
#    instantanée sentence, and drop the stale spell-check markers
#    that surrounded the misspelled word.
# ------------------------------------------------------------------
$find = $d.Content
$ok = $find.Find.Execute("messagerie instantanné.", $true, $false, $false,
                          $false, $false, $true, 1, $false, $null, 0)

if ($ok) {
    $editStart = $find.Start
    $newSentence = "messagerie instantané."
    $find.Text = $newSentence

    # Position right before the trailing "." - this is where Word
    # drops the _GoBack bookmark after the last edit made in the doc.
    $bmPos = $editStart + $newSentence.Length - 1
    $bmRange = $d.Range($bmPos, $bmPos)

    # ------------------------------------------------------------------
    # 2) Move the "_GoBack" bookmark to the location of this edit. Word
    #    only ever keeps a single "_GoBack" bookmark in a document, so
    #    adding it here automatically removes the old one that used to
    #    sit after the "Cette fonctionnalité ajoutera ..." paragraph.
    # ------------------------------------------------------------------
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
